# edit.ps1 - applies the commit:
#   1. Bumps the cached datetimeFigureOut field text 25/03/2017 -> 26/03/2017
#      (NotesMaster + HandoutMaster "Date Placeholder" shapes).
#   2. Changes the cached slidenum field glyph from the anonymous numeral
#      mark <#> to <nr.> on the SlideMaster, all 7 CustomLayouts, the
#      NotesMaster and the HandoutMaster ("Slide Number Placeholder" shapes).
#   3. On slide 16, splits the "2 webservers en aparte storage server" bullet
#      into three runs, inserting a red "(/fileserver?)" remark after
#      "2 webservers".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: find the first shape in a Shapes collection with a given
# PlaceholderFormat.Type (13 = slide number, 16 = date).
# ---------------------------------------------------------------------------
function Get-PlaceholderShape($shapes, $phType) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $t = $null
        try { $t = $sh.PlaceholderFormat.Type } catch { $t = $null }
        if ($t -eq $phType) {
            return $sh
        }
    }
    return $null
}

$ppPlaceholderDate = 16
$ppPlaceholderSlideNumber = 13

# ---------------------------------------------------------------------------
# 1 & 2. SlideMaster: slide-number glyph.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

$mSlideNum = Get-PlaceholderShape $master.Shapes $ppPlaceholderSlideNumber
if ($mSlideNum -ne $null) {
    $mSlideNum.TextFrame.TextRange.Text = "‹nr.›"
}

# ---------------------------------------------------------------------------
# 2. All seven CustomLayouts: slide-number glyph.
# ---------------------------------------------------------------------------
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $lSlideNum = Get-PlaceholderShape $layout.Shapes $ppPlaceholderSlideNumber
    if ($lSlideNum -ne $null) {
        $lSlideNum.TextFrame.TextRange.Text = "‹nr.›"
    }
}

# ---------------------------------------------------------------------------
# 1 & 2. NotesMaster: date + slide-number glyph.
# ---------------------------------------------------------------------------
$notesMaster = $p.NotesMaster

$nmDate = Get-PlaceholderShape $notesMaster.Shapes $ppPlaceholderDate
if ($nmDate -ne $null) {
    $nmDate.TextFrame.TextRange.Text = "26/03/2017"
}

$nmSlideNum = Get-PlaceholderShape $notesMaster.Shapes $ppPlaceholderSlideNumber
if ($nmSlideNum -ne $null) {
    $nmSlideNum.TextFrame.TextRange.Text = "‹nr.›"
}

# ---------------------------------------------------------------------------
# 1 & 2. HandoutMaster: date + slide-number glyph.
# ---------------------------------------------------------------------------
$handoutMaster = $p.HandoutMaster

$hmDate = Get-PlaceholderShape $handoutMaster.Shapes $ppPlaceholderDate
if ($hmDate -ne $null) {
    $hmDate.TextFrame.TextRange.Text = "26/03/2017"
}

$hmSlideNum = Get-PlaceholderShape $handoutMaster.Shapes $ppPlaceholderSlideNumber
if ($hmSlideNum -ne $null) {
    $hmSlideNum.TextFrame.TextRange.Text = "‹nr.›"
}

# ---------------------------------------------------------------------------
# 3. Slide 16: split the "2 webservers en aparte storage server" bullet and
#    insert a red "(/fileserver?)" remark.
# ---------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
$contentShape = $slide16.Shapes.Item(1)
$bodyRange = $contentShape.TextFrame.TextRange

$targetParaIndex = -1
for ($j = 1; $j -le $bodyRange.Paragraphs().Count; $j++) {
    $paraText = $bodyRange.Paragraphs($j, 1).Text.TrimEnd("`r")
    if ($paraText -eq "2 webservers en aparte storage server") {
        $targetParaIndex = $j
        break
    }
}

if ($targetParaIndex -ne -1) {
    $para = $bodyRange.Paragraphs($targetParaIndex, 1)

    # Retype the whole line with the new wording inserted.
    $para.Text = "2 webservers(/fileserver?) en aparte storage server"

    # Re-fetch the paragraph and colour just the inserted remark red.
    $para = $bodyRange.Paragraphs($targetParaIndex, 1)
    $remark = $para.Characters(13, 14)
    $remark.Font.Color.RGB = 255
}
